$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 133-162 per diff ---
# Row 133
$ws.Cells.Item(133, 4).Value = 44559
$ws.Cells.Item(133, 10).Value = 300
$ws.Cells.Item(133, 11).Value = 14000
$ws.Cells.Item(133, 12).Value = 15000
$ws.Cells.Item(133, 13).Value = 14500
$ws.Cells.Item(133, 14).Value = '$/caja 12 unidades'
$ws.Cells.Item(133, 16).Value = 1208
$ws.Cells.Item(133, 17).Value = 12

# Row 134
$ws.Cells.Item(134, 4).Value = 44559
$ws.Cells.Item(134, 8).Value = 'Tuna'
$ws.Cells.Item(134, 9).Value = 'Extra'
$ws.Cells.Item(134, 10).Value = 300
$ws.Cells.Item(134, 11).Value = 14000
$ws.Cells.Item(134, 12).Value = 15000
$ws.Cells.Item(134, 13).Value = 14500
$ws.Cells.Item(134, 14).Value = '$/caja 12 unidades'
$ws.Cells.Item(134, 16).Value = 1208
$ws.Cells.Item(134, 17).Value = 12

# Row 135
$ws.Cells.Item(135, 8).Value = 'Calameño'
$ws.Cells.Item(135, 9).Value = 'Primera'
$ws.Cells.Item(135, 11).Value = 1000
$ws.Cells.Item(135, 12).Value = 1200
$ws.Cells.Item(135, 13).Value = 1100
$ws.Cells.Item(135, 16).Value = 1100

# Row 136
$ws.Cells.Item(136, 8).Value = 'Calameño'

# Row 137
$ws.Cells.Item(137, 4).Value = 44208
$ws.Cells.Item(137, 8).Value = 'Tuna'
$ws.Cells.Item(137, 9).Value = 'Extra'
$ws.Cells.Item(137, 11).Value = 1400
$ws.Cells.Item(137, 12).Value = 1400
$ws.Cells.Item(137, 13).Value = 1400
$ws.Cells.Item(137, 16).Value = 1400

# Row 138
$ws.Cells.Item(138, 4).Value = 44208
$ws.Cells.Item(138, 10).Value = 2000
$ws.Cells.Item(138, 13).Value = 1100
$ws.Cells.Item(138, 16).Value = 1100

# Row 139
$ws.Cells.Item(139, 4).Value = 44264
$ws.Cells.Item(139, 9).Value = 'Primera'
$ws.Cells.Item(139, 10).Value = 2000
$ws.Cells.Item(139, 11).Value = 1000
$ws.Cells.Item(139, 12).Value = 1200
$ws.Cells.Item(139, 13).Value = 1100
$ws.Cells.Item(139, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(139, 16).Value = 1100

# Row 140
$ws.Cells.Item(140, 4).Value = 44264
$ws.Cells.Item(140, 8).Value = 'Tuna'
$ws.Cells.Item(140, 10).Value = 2500
$ws.Cells.Item(140, 13).Value = 1120
$ws.Cells.Item(140, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(140, 16).Value = 1120

# Row 141
$ws.Cells.Item(141, 8).Value = 'Calameño'

# Row 142
$ws.Cells.Item(142, 8).Value = 'Calameño'

# Row 143
$ws.Cells.Item(143, 4).Value = 44232
$ws.Cells.Item(143, 8).Value = 'Tuna'
$ws.Cells.Item(143, 10).Value = 4000
$ws.Cells.Item(143, 11).Value = 1400
$ws.Cells.Item(143, 12).Value = 1500
$ws.Cells.Item(143, 13).Value = 1450
$ws.Cells.Item(143, 14).Value = '$/unidad'
$ws.Cells.Item(143, 15).Value = 'Región del Maule'
$ws.Cells.Item(143, 16).Value = 1450
$ws.Cells.Item(143, 17).Value = 1

# Row 144
$ws.Cells.Item(144, 4).Value = 44232
$ws.Cells.Item(144, 8).Value = 'Tuna'
$ws.Cells.Item(144, 10).Value = 4000
$ws.Cells.Item(144, 11).Value = 1000
$ws.Cells.Item(144, 12).Value = 1200
$ws.Cells.Item(144, 13).Value = 1100
$ws.Cells.Item(144, 15).Value = 'Región del Maule'
$ws.Cells.Item(144, 16).Value = 1100

# Row 145
$ws.Cells.Item(145, 9).Value = 'Extra'
$ws.Cells.Item(145, 10).Value = 300
$ws.Cells.Item(145, 11).Value = 15000
$ws.Cells.Item(145, 12).Value = 15000
$ws.Cells.Item(145, 13).Value = 15000
$ws.Cells.Item(145, 14).Value = '$/caja 12 unidades'
$ws.Cells.Item(145, 16).Value = 1250
$ws.Cells.Item(145, 17).Value = 12

# Row 146
$ws.Cells.Item(146, 8).Value = 'Calameño'
$ws.Cells.Item(146, 9).Value = 'Primera'
$ws.Cells.Item(146, 10).Value = 1500
$ws.Cells.Item(146, 11).Value = 1500
$ws.Cells.Item(146, 12).Value = 1500
$ws.Cells.Item(146, 13).Value = 1500
$ws.Cells.Item(146, 14).Value = '$/unidad'
$ws.Cells.Item(146, 16).Value = 1500
$ws.Cells.Item(146, 17).Value = 1

# Row 147
$ws.Cells.Item(147, 8).Value = 'Calameño'
$ws.Cells.Item(147, 9).Value = 'Segunda'
$ws.Cells.Item(147, 11).Value = 1000
$ws.Cells.Item(147, 12).Value = 1000
$ws.Cells.Item(147, 13).Value = 1000
$ws.Cells.Item(147, 16).Value = 1000

# Row 148
$ws.Cells.Item(148, 9).Value = 'Extra'
$ws.Cells.Item(148, 10).Value = 300
$ws.Cells.Item(148, 11).Value = 14000
$ws.Cells.Item(148, 12).Value = 14000
$ws.Cells.Item(148, 13).Value = 14000
$ws.Cells.Item(148, 14).Value = '$/caja 12 unidades'
$ws.Cells.Item(148, 16).Value = 1167
$ws.Cells.Item(148, 17).Value = 12

# Row 149
$ws.Cells.Item(149, 4).Value = 44551
$ws.Cells.Item(149, 8).Value = 'Tuna'
$ws.Cells.Item(149, 9).Value = 'Primera'
$ws.Cells.Item(149, 10).Value = 1500
$ws.Cells.Item(149, 15).Value = 'Región de O''Higgins'

# Row 150
$ws.Cells.Item(150, 4).Value = 44551
$ws.Cells.Item(150, 8).Value = 'Tuna'
$ws.Cells.Item(150, 9).Value = 'Segunda'
$ws.Cells.Item(150, 10).Value = 1500
$ws.Cells.Item(150, 11).Value = 1000
$ws.Cells.Item(150, 12).Value = 1000
$ws.Cells.Item(150, 13).Value = 1000
$ws.Cells.Item(150, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(150, 16).Value = 1000

# Row 151
$ws.Cells.Item(151, 8).Value = 'Calameño'
$ws.Cells.Item(151, 15).Value = 'Región de O''Higgins'

# Row 152
$ws.Cells.Item(152, 8).Value = 'Calameño'
$ws.Cells.Item(152, 15).Value = 'Región de O''Higgins'

# Row 153
$ws.Cells.Item(153, 4).Value = 44238
$ws.Cells.Item(153, 8).Value = 'Tuna'
$ws.Cells.Item(153, 11).Value = 1500
$ws.Cells.Item(153, 12).Value = 1500
$ws.Cells.Item(153, 13).Value = 1500
$ws.Cells.Item(153, 15).Value = 'Región del Maule'
$ws.Cells.Item(153, 16).Value = 1500

# Row 154
$ws.Cells.Item(154, 4).Value = 44238
$ws.Cells.Item(154, 9).Value = 'Primera'
$ws.Cells.Item(154, 11).Value = 1200
$ws.Cells.Item(154, 12).Value = 1200
$ws.Cells.Item(154, 13).Value = 1200
$ws.Cells.Item(154, 15).Value = 'Región del Maule'
$ws.Cells.Item(154, 16).Value = 1200

# Row 155
$ws.Cells.Item(155, 4).Value = 44257
$ws.Cells.Item(155, 9).Value = 'Extra'
$ws.Cells.Item(155, 10).Value = 3000
$ws.Cells.Item(155, 11).Value = 1400
$ws.Cells.Item(155, 12).Value = 1400
$ws.Cells.Item(155, 13).Value = 1400
$ws.Cells.Item(155, 16).Value = 1400

# Row 156
$ws.Cells.Item(156, 4).Value = 44257
$ws.Cells.Item(156, 9).Value = 'Extra'
$ws.Cells.Item(156, 10).Value = 3000
$ws.Cells.Item(156, 11).Value = 1400
$ws.Cells.Item(156, 12).Value = 1400
$ws.Cells.Item(156, 13).Value = 1400
$ws.Cells.Item(156, 16).Value = 1400

# Row 157
$ws.Cells.Item(157, 4).Value = 44200
$ws.Cells.Item(157, 10).Value = 1500
$ws.Cells.Item(157, 12).Value = 1000
$ws.Cells.Item(157, 13).Value = 1000
$ws.Cells.Item(157, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(157, 16).Value = 1000

# Row 158
$ws.Cells.Item(158, 4).Value = 44200
$ws.Cells.Item(158, 10).Value = 1500
$ws.Cells.Item(158, 12).Value = 1000
$ws.Cells.Item(158, 13).Value = 1000
$ws.Cells.Item(158, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(158, 16).Value = 1000

# Row 159
$ws.Cells.Item(159, 4).Value = 44224
$ws.Cells.Item(159, 9).Value = 'Primera'
$ws.Cells.Item(159, 10).Value = 2000
$ws.Cells.Item(159, 11).Value = 1000
$ws.Cells.Item(159, 12).Value = 1200
$ws.Cells.Item(159, 13).Value = 1100
$ws.Cells.Item(159, 16).Value = 1100

# Row 160
$ws.Cells.Item(160, 4).Value = 44224
$ws.Cells.Item(160, 8).Value = 'Tuna'
$ws.Cells.Item(160, 10).Value = 2000

# Row 161
$ws.Cells.Item(161, 8).Value = 'Calameño'

# Row 162
$ws.Cells.Item(162, 8).Value = 'Calameño'

# --- Add new rows 163-164 ---
# Row 163
$ws.Cells.Item(163, 1).Value = 4
$ws.Cells.Item(163, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(163, 3).Value = 'Los Lagos'
$ws.Cells.Item(163, 4).Value = 44239
$ws.Cells.Item(163, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(163, 5).Value = 10
$ws.Cells.Item(163, 6).Value = 100112027
$ws.Cells.Item(163, 7).Value = 'Melón'
$ws.Cells.Item(163, 8).Value = 'Tuna'
$ws.Cells.Item(163, 9).Value = 'Extra'
$ws.Cells.Item(163, 10).Value = 3000
$ws.Cells.Item(163, 11).Value = 1400
$ws.Cells.Item(163, 12).Value = 1500
$ws.Cells.Item(163, 13).Value = 1450
$ws.Cells.Item(163, 14).Value = '$/unidad'
$ws.Cells.Item(163, 15).Value = 'Región del Maule'
$ws.Cells.Item(163, 16).Value = 1450
$ws.Cells.Item(163, 17).Value = 1
$ws.Cells.Item(163, 18).Value = 'Hortaliza'

# Row 164
$ws.Cells.Item(164, 1).Value = 4
$ws.Cells.Item(164, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(164, 3).Value = 'Los Lagos'
$ws.Cells.Item(164, 4).Value = 44239
$ws.Cells.Item(164, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(164, 5).Value = 10
$ws.Cells.Item(164, 6).Value = 100112027
$ws.Cells.Item(164, 7).Value = 'Melón'
$ws.Cells.Item(164, 8).Value = 'Tuna'
$ws.Cells.Item(164, 9).Value = 'Primera'
$ws.Cells.Item(164, 10).Value = 3000
$ws.Cells.Item(164, 11).Value = 1000
$ws.Cells.Item(164, 12).Value = 1200
$ws.Cells.Item(164, 13).Value = 1100
$ws.Cells.Item(164, 14).Value = '$/unidad'
$ws.Cells.Item(164, 15).Value = 'Región del Maule'
$ws.Cells.Item(164, 16).Value = 1100
$ws.Cells.Item(164, 17).Value = 1
$ws.Cells.Item(164, 18).Value = 'Hortaliza'
